$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff swaps values between row 2 and row 3 for columns D, M, N, P, S
# (O stays the same at 6000 for both rows).

# Column D (Fecha)
$ws.Range("D2").Value = 44991
$ws.Range("D3").Value = 44995

# Column M (Volumen)
$ws.Range("M2").Value = 50
$ws.Range("M3").Value = 100

# Column N (Precio minimo)
$ws.Range("N2").Value = 6000
$ws.Range("N3").Value = 5500

# Column P (Precio promedio ponderado)
$ws.Range("P2").Value = 6000
$ws.Range("P3").Value = 5750

# Column S (Precio $/Kg)
$ws.Range("S2").Value = 3000
$ws.Range("S3").Value = 2875
